$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4189.8423
$ws.Range("I116").Value = 2061
$ws.Range("J116").Value = 4950.143
$ws.Range("K116").Value = 2061
$ws.Range("L116").Value = 4950.143
$ws.Range("M116").Value = 1381
$ws.Range("N116").Value = -11834.143
$ws.Range("H129").Value = 295344.34
$ws.Range("J129").Value = 323898.66
$ws.Range("L129").Value = 971695.98
$ws.Range("N129").Value = -981695.98
$ws.Range("H137").Value = 2590.15
$ws.Range("I137").Value = 3063.9092
$ws.Range("J137").Value = 2011.1111
$ws.Range("K137").Value = 9191.7276
$ws.Range("L137").Value = 6033.3333
$ws.Range("M137").Value = -6641.7276
$ws.Range("N137").Value = -11133.3333
$ws.Range("H138").Value = 1904.202
$ws.Range("I138").Value = 372.72342
$ws.Range("J138").Value = 3288.423
$ws.Range("K138").Value = 1118.17026
$ws.Range("L138").Value = 9865.269
$ws.Range("M138").Value = 4021.82974
$ws.Range("N138").Value = -20145.269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 462655.7
$ws.Range("I61").Value = 501127.16
$ws.Range("J61").Value = 998
$ws.Range("K61").Value = 501127.16
$ws.Range("L61").Value = 998
$ws.Range("M61").Value = -500915.16
$ws.Range("N61").Value = -1422
$ws.Range("H80").Value = 65535
$ws.Range("J80").Value = 65535
$ws.Range("L80").Value = 65535
$ws.Range("N80").Value = -67531
$ws.Range("H82").Value = 15000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 65535
$ws.Range("J83").Value = 65535
$ws.Range("L83").Value = 196605
$ws.Range("N83").Value = -206589
$ws.Range("H85").Value = 15000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H102").Value = 1366.6666
$ws.Range("I102").Value = 1050
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1050
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 572
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 17013.758
$ws.Range("I132").Value = 1871.8518
$ws.Range("K132").Value = 5615.555399999999
$ws.Range("M132").Value = -3085.555399999999
$ws.Range("H136").Value = 462655.7
$ws.Range("I136").Value = 501127.16
$ws.Range("J136").Value = 998
$ws.Range("K136").Value = 1503381.48
$ws.Range("L136").Value = 2994
$ws.Range("M136").Value = -1500831.48
$ws.Range("N136").Value = -8094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 9200
$ws.Range("J76").Value = 9200
$ws.Range("L76").Value = 9200
$ws.Range("N76").Value = -9830
$ws.Range("H79").Value = 9200
$ws.Range("J79").Value = 9200
$ws.Range("L79").Value = 9200
$ws.Range("N79").Value = -11384
$ws.Range("H99").Value = 1777.1818
$ws.Range("I99").Value = 1962.25
$ws.Range("K99").Value = 1962.25
$ws.Range("M99").Value = -464.25
$ws.Range("H134").Value = 7160.5293
$ws.Range("I134").Value = 9448.091
$ws.Range("K134").Value = 28344.273
$ws.Range("M134").Value = -25809.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 552
$ws.Range("H31").Value = 4505.731
$ws.Range("I31").Value = 2950.7693
$ws.Range("K31").Value = 2950.7693
$ws.Range("M31").Value = -2655.7693
$ws.Range("H34").Value = 4505.731
$ws.Range("I34").Value = 2950.7693
$ws.Range("K34").Value = 2950.7693
$ws.Range("M34").Value = -2748.7693
$ws.Range("H62").Value = 35717500
$ws.Range("I62").Value = 40002796
$ws.Range("J62").Value = 6668.6665
$ws.Range("K62").Value = 40002796
$ws.Range("L62").Value = 6668.6665
$ws.Range("M62").Value = -40002172
$ws.Range("N62").Value = -7916.6665
$ws.Range("H65").Value = 35717500
$ws.Range("I65").Value = 40002796
$ws.Range("J65").Value = 6668.6665
$ws.Range("K65").Value = 200013980
$ws.Range("L65").Value = 33343.3325
$ws.Range("M65").Value = -200010860
$ws.Range("N65").Value = -39583.3325
$ws.Range("H103").Value = 18500
$ws.Range("I103").Value = 18500
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 18500
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -17328
$ws.Range("N103").ClearContents()
$ws.Range("H122").Value = 2368.3333
$ws.Range("I122").Value = 2557.2856
$ws.Range("J122").Value = 1707
$ws.Range("K122").Value = 7671.8568
$ws.Range("L122").Value = 5121
$ws.Range("M122").Value = -5221.8568
$ws.Range("N122").Value = -10021
$ws.Range("H132").Value = 1788.5106
$ws.Range("I132").Value = 1441.8096
$ws.Range("J132").Value = 4700.8
$ws.Range("K132").Value = 4325.4288
$ws.Range("L132").Value = 14102.4
$ws.Range("M132").Value = -1795.4288
$ws.Range("N132").Value = -19162.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 707.0303
$ws.Range("J131").Value = 723.3626
$ws.Range("L131").Value = 2170.0878
$ws.Range("N131").Value = -12250.0878
$ws.Range("H134").Value = 1912.84
$ws.Range("I134").Value = 1421.8
$ws.Range("J134").Value = 3877
$ws.Range("K134").Value = 4265.4
$ws.Range("L134").Value = 11631
$ws.Range("M134").Value = 804.6000000000004
$ws.Range("N134").Value = -21771

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 93129.836
$ws.Range("I132").Value = 11755.8
$ws.Range("J132").Value = 500000
$ws.Range("K132").Value = 35267.39999999999
$ws.Range("L132").Value = 1500000
$ws.Range("M132").Value = -32737.39999999999
$ws.Range("N132").Value = -1505060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3285.4736
$ws.Range("I40").Value = 2995.6428
$ws.Range("K40").Value = 2995.6428
$ws.Range("M40").Value = -2859.6428
$ws.Range("H68").Value = 2228
$ws.Range("I68").Value = 1649.75
$ws.Range("K68").Value = 1649.75
$ws.Range("M68").Value = -900.75
$ws.Range("H71").Value = 2228
$ws.Range("I71").Value = 1649.75
$ws.Range("K71").Value = 8248.75
$ws.Range("M71").Value = -4504.75
$ws.Range("H122").Value = 855373.5
$ws.Range("I122").Value = 1963429.2
$ws.Range("J122").Value = 3023
$ws.Range("K122").Value = 5890287.6
$ws.Range("L122").Value = 9069
$ws.Range("M122").Value = -5887837.6
$ws.Range("N122").Value = -13969
$ws.Range("H132").Value = 635726.0600000001
$ws.Range("I132").Value = 635726.0600000001
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1907178.18
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1904648.18
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1087.3334
$ws.Range("I136").Value = 952.23334
$ws.Range("K136").Value = 2856.70002
$ws.Range("M136").Value = -306.7000200000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1476.3529
$ws.Range("I122").Value = 1476.3529
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4429.0587
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1979.0587
$ws.Range("N122").ClearContents()
